$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 10 (Objetivos / Objectives Portuguese text was missing and
#     had the professor name in its place by mistake) ---
$objetivosPt = "Fornecer aos alunos capacidade de elaborar projetos de sistemas de tratamento de águas residuárias, envolvendo: memorial justificativa e de cálculos, elaborar manual de operação, especificar equipamentos e elaborar planos de controle operacional e de emergência."
$ws.Cells.Item(10,2).Value = $objetivosPt
$ws.Cells.Item(10,3).Value = $objetivosPt

# --- Insert a new row at 13 for "Docentes responsáveis:" content,
#     shifting the existing rows 13-21 down to 14-22 ---
$ws.Rows.Item(13).Insert()

# The insert leaves a stray formatted-but-empty cell in column A; remove it
# entirely so row 13 only carries the B/C content (matches target layout).
$ws.Cells.Item(13,1).Clear()

# Copy formatting from the row above (B11/C11, s=2/s=3) onto the new B13/C13
# cells so we don't leave unused/duplicate styles behind.
$ws.Cells.Item(11,2).Copy() | Out-Null
$ws.Cells.Item(13,2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(11,3).Copy() | Out-Null
$ws.Cells.Item(13,3).PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(13,2).Value = "7455355 - Robson da Silva Rocha"
$ws.Cells.Item(13,3).Value = "7455355 - Robson da Silva Rocha"

# --- Row 14 (was 13): "Programa resumido:" short syllabus text ---
$programaResumido = "Métodos alternativos para o tratamento de áquas residuárias e legislação."
$ws.Cells.Item(14,2).Value = $programaResumido
$ws.Cells.Item(14,3).Value = $programaResumido

# --- Row 16 (was 15): "Programa:" full syllabus text ---
$programa = "Revisão sobre alternativas para tratamento de águas residuárias e legislação sobre padrões de emissão e de qualidade. Formas de apresentação e conteúdo de Relatórios de Avaliação Preliminar (RAP) e de projetos de sistemas de águas residuárias, com base em normas da ABNT e de órgãos de controlede poluição. Equipamentos eletromecânicos e eletrônicos mais utilizados em sistemas de tratamento de águas residuárias (STAR): especificação e manutenção. Informatização e automação de STARs. Paisagísmo de STARs. Gerenciamento de STARs."
$ws.Cells.Item(16,2).Value = $programa
$ws.Cells.Item(16,3).Value = $programa

# --- Row 19 (was 18): "Método:" text ---
$metodo = 'Serão ministradas aulas expositivas convencionais, associadas a apresentação de vídeos e "slides"sobre sistemas de tratamento, visitas a instalações existentes. Além disso os alunos deverão acompanhar a operação (e elaborar relatórios) de estações de tratamento em funcionamento.'
$ws.Cells.Item(19,2).Value = $metodo
$ws.Cells.Item(19,3).Value = $metodo

# --- Row 20 (was 19): "Critério:" text ---
$criterio = "Os alunos serão submetidos a duas provas escritas, sem consulta (1 e 2 bimestres)."
$ws.Cells.Item(20,2).Value = $criterio
$ws.Cells.Item(20,3).Value = $criterio

# --- Row 21 (was 20): "Norma de recuperação:" text ---
$normaRecuperacao = "Elaboração de monografia, com tema escolhido pelo docente, enfocando matéria em que o aluno demonstrou menor habilifdade (peso: 3,0); e prova escrita sobre toda a matéria da disciplina (peso: 7)."
$ws.Cells.Item(21,2).Value = $normaRecuperacao
$ws.Cells.Item(21,3).Value = $normaRecuperacao

# --- Row 22 (was 21): "Bibliografia:" text ---
$bibliografia = "PAGANINI, W.S. (1997). Disposição de Esgoto no Solo, AESABESP, 2a Ed., 232p.`nKELLNER, E.; PIRES, E.C. (1996). Lagoas de Estabilização; projeto e Operação. Rio de Janeiro (RJ). Brasil, Ed. ABES, 241 p. `nVAN HAANDEL, A.; MARAIS, G. (1999). O Comportamento do Sistema de Lodo Ativado: Teoria e Operações para Projeto e Operação. Universidade Federal da Paraíba - epgraf - Campina Grande, Pb, 477p. `nCatálogos de fornecedores de materiais e de equipamentos utilizados em sistemas de tratamento de águas residuárias."
$ws.Cells.Item(22,2).Value = $bibliografia
$ws.Cells.Item(22,3).Value = $bibliografia
